$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Suggestions")

# Fill in the new Suggestions data
$ws3.Range("A2").Value = "Tinder"
$ws3.Range("B2").Value = "Mad Paws"
$ws3.Range("C2").Value = "Hoyts Group"
$ws3.Range("D2").Value = "Event Cinema"
$ws3.Range("E2").Value = "Gather Online"
$ws3.Range("F2").Value = "Clique Labs"

$ws3.Range("A3").Value = "Clique Labs"
$ws3.Range("B3").Value = "Insitchu"
$ws3.Range("C3").Value = "Glamcorner"
$ws3.Range("D3").Value = "Deliveroo"
$ws3.Range("E3").Value = "Happn"
$ws3.Range("F3").Value = "The Iconic"
$ws3.Range("G3").Value = "Coles"
$ws3.Range("H3").Value = "Canva"
$ws3.Range("I3").Value = "DesignCrowd"
$ws3.Range("J3").Value = "Mastercard"

# Apply formatting across the used range
$ws3.Range("A1:K5").Font.Color = 0
$ws3.Range("A2").NumberFormat = "@"
$ws3.Range("A4").NumberFormat = "@"

# Update the selection/view and active sheet
$ws3.Activate()
$ws3.Range("A1:K3").Select() | Out-Null
